$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update column B for existing rows 329-385 (previously placeholder 0 values) ---
$updateValues = @(455,454,458,459,463,454,451,456,480,557,555,563,569,645,644,640,648,728,758,753,756,805,817,819,824,817,842,841,842,782,801,825,845,748,746,745,754,787,801,771,758,726,717,714,704,604,584,550,542,425,422,298,278,203,190,176,171)
$row = 329
foreach ($v in $updateValues) {
    $ws.Cells.Item($row, 2).Value = $v
    $row = $row + 1
}

# --- Apply the date/time number format used by column A to the new rows before setting values ---
$ws.Range("A386:A481").NumberFormat = $ws.Range("A385").NumberFormat

# --- Append new rows 386-481 with timestamps (column A) and production values (column B) ---
$newTimestamps = @(45734,45734.01041666666,45734.02083333334,45734.03125,45734.04166666666,45734.05208333334,45734.0625,45734.07291666666,45734.08333333334,45734.09375,45734.10416666666,45734.11458333334,45734.125,45734.13541666666,45734.14583333334,45734.15625,45734.16666666666,45734.17708333334,45734.1875,45734.19791666666,45734.20833333334,45734.21875,45734.22916666666,45734.23958333334,45734.25,45734.26041666666,45734.27083333334,45734.28125,45734.29166666666,45734.30208333334,45734.3125,45734.32291666666,45734.33333333334,45734.34375,45734.35416666666,45734.36458333334,45734.375,45734.38541666666,45734.39583333334,45734.40625,45734.41666666666,45734.42708333334,45734.4375,45734.44791666666,45734.45833333334,45734.46875,45734.47916666666,45734.48958333334,45734.5,45734.51041666666,45734.52083333334,45734.53125,45734.54166666666,45734.55208333334,45734.5625,45734.57291666666,45734.58333333334,45734.59375,45734.60416666666,45734.61458333334,45734.625,45734.63541666666,45734.64583333334,45734.65625,45734.66666666666,45734.67708333334,45734.6875,45734.69791666666,45734.70833333334,45734.71875,45734.72916666666,45734.73958333334,45734.75,45734.76041666666,45734.77083333334,45734.78125,45734.79166666666,45734.80208333334,45734.8125,45734.82291666666,45734.83333333334,45734.84375,45734.85416666666,45734.86458333334,45734.875,45734.88541666666,45734.89583333334,45734.90625,45734.91666666666,45734.92708333334,45734.9375,45734.94791666666,45734.95833333334,45734.96875,45734.97916666666,45734.98958333334)
$newValues = @(212,195,195,194,191,195,222,227,229,193,190,178,173,185,197,179,211,222,232,256,232,240,247,277,293,354,340,328,359,337,327,334,244,229,228,227,156,143,137,135,128,121,105,97,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0)
$row = 386
for ($i = 0; $i -lt $newTimestamps.Length; $i++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamps[$i]
    $ws.Cells.Item($row, 2).Value = $newValues[$i]
    $row = $row + 1
}
